$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Row 10: PO line updated (new PO number, new delivery dates, qty)
# ---------------------------------------------------------------
$ws.Range("A10").Value = 151423149
$ws.Range("E10").Value = 45308
$ws.Range("J10").Value = 10
$ws.Range("K10").Value = 45321

# ---------------------------------------------------------------
# Row 18: PO line updated
# ---------------------------------------------------------------
$ws.Range("A18").Value = 151423149
$ws.Range("C18").Value = 20
$ws.Range("E18").Value = 45308
$ws.Range("J18").Value = 10
$ws.Range("K18").Value = 45321

# ---------------------------------------------------------------
# Row 34: line replaced with a new APS PO line (new formatting too)
# ---------------------------------------------------------------
$ws.Range("A71:K71").Copy()
$ws.Range("A34:K34").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A34").Value = 151448381
$ws.Range("B34").Value = "APS"
$ws.Range("C34").Value = 40
$ws.Range("D34").Value = 3555
$ws.Range("E34").Value = 45362
$ws.Range("F34").Value = 646872
$ws.Range("G34").Value = "Water Bellow_BP60_Dia 170x300mm"
$ws.Range("H34").Value = 10
$ws.Range("I34").Value = "PCS"
$ws.Range("J34").Value = 10
$ws.Range("K34").Value = 45392

# ---------------------------------------------------------------
# Row 40: PO number / qty / dates updated (formatting unchanged)
# ---------------------------------------------------------------
$ws.Range("A40").Value = 151448381
$ws.Range("C40").Value = 50
$ws.Range("E40").Value = 45362
$ws.Range("H40").Value = 5
$ws.Range("J40").Value = 5
$ws.Range("K40").Value = 45392

# ---------------------------------------------------------------
# Row 42: line replaced with a new APS PO line (new formatting too)
# ---------------------------------------------------------------
$ws.Range("A40:K40").Copy()
$ws.Range("A42:K42").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("K48").Copy()
$ws.Range("K42").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A42").Value = 151448381
$ws.Range("B42").Value = "APS"
$ws.Range("C42").Value = 20
$ws.Range("D42").Value = 5130
$ws.Range("E42").Value = 45362
$ws.Range("F42").Value = 663092
$ws.Range("G42").Value = "Cement Bellow_Dia330x200 lg_MT3.0"
$ws.Range("H42").Value = 10
$ws.Range("I42").Value = "PCS"
$ws.Range("J42").Value = 10
$ws.Range("K42").Value = 45392

# ---------------------------------------------------------------
# Row 44: line replaced with a new APS PO line (new formatting too)
# ---------------------------------------------------------------
$ws.Range("A40:K40").Copy()
$ws.Range("A44:K44").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("K48").Copy()
$ws.Range("K44").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A44").Value = 151448381
$ws.Range("B44").Value = "APS"
$ws.Range("C44").Value = 30
$ws.Range("D44").Value = 6615
$ws.Range("E44").Value = 45362
$ws.Range("F44").Value = 663093
$ws.Range("G44").Value = "Screw conv. Bellow_Dia278 x500lg_MT3.0"
$ws.Range("H44").Value = 10
$ws.Range("I44").Value = "PCS"
$ws.Range("J44").Value = 10
$ws.Range("K44").Value = 45390

# ---------------------------------------------------------------
# Rows 51 / 55 / 62: blank spacer cell restyled (G column)
# ---------------------------------------------------------------
$ws.Range("G25").Copy()
$ws.Range("G51").PasteSpecial(-4122)
$ws.Range("G55").PasteSpecial(-4122)
$ws.Range("G62").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# Insert a brand-new PO line as row 68 (everything below shifts
# down by one row, through the former row 79 which becomes row 80)
# ---------------------------------------------------------------
$ws.Rows("68:68").Insert()

$ws.Range("A6:K6").Copy()
$ws.Range("A68:K68").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("K48").Copy()
$ws.Range("K68").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A68").Value = 151448381
$ws.Range("B68").Value = "APS"
$ws.Range("C68").Value = 10
$ws.Range("D68").Value = 1896
$ws.Range("E68").Value = 45362
$ws.Range("F68").Value = 716663
$ws.Range("G68").Value = "BELLOW F. water hopper MT 1.0 IBC"
$ws.Range("H68").Value = 2
$ws.Range("I68").Value = "PCS"
$ws.Range("J68").Value = 2
$ws.Range("K68").Value = 45383

# ---------------------------------------------------------------
# Final view state: scroll back to top, select C16
# ---------------------------------------------------------------
$ws.Range("C16").Select()
